$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.199.84"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.903.08"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'308.21"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.5210"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").Value = "'0.3765"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'0.07278"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "'0.9054"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'0.08253"
$ws.Range("E12").Value = "  +7.89%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'96.88"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.900.37"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "'5.294"
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'0.000008680"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "27.233.28"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'5.097"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").Value = "2.161.56"
$ws.Range("E22").Value = "  +2.94%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'6.437"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'2.327"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "'146.44"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'1.748"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").Value = "'18.23"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'114.99"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").Value = "'4.835"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "'4.904"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'0.09269"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "'0.05085"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "'0.7998"
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'3.426"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").Value = "'2.947"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "'2.597"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'0.5717"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'6.594"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'117.03"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "'0.4857"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("D48").Value = "'10.11"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "'1.630"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "'64.01"
$ws.Range("E51").Value = "  +0.45%  "
